$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 71, pushing the existing rows 71-80 down to 74-83.
$ws.Range("A71:T73").EntireRow.Insert()

# New row 71: Primera, week of 2023-08-03 (serial 45141)
$ws.Range("A71").Value = 1
$ws.Range("B71").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C71").Value = "Arica y Parinacota"
$ws.Range("D71").Value = 45141
$ws.Range("E71").Value = 15
$ws.Range("F71").Value = "Fruta"
$ws.Range("G71").Value = 100101
$ws.Range("H71").Value = "Berries"
$ws.Range("I71").Value = 100112025
$ws.Range("J71").Value = "Frutilla"
$ws.Range("K71").Value = "Sin especificar"
$ws.Range("L71").Value = "Primera"
$ws.Range("M71").Value = 70
$ws.Range("N71").Value = 7000
$ws.Range("O71").Value = 8000
$ws.Range("P71").Value = 7500
$ws.Range("Q71").Value = "$/bandeja 3 kilos"
$ws.Range("R71").Value = "Región de Arica y Parinacota"
$ws.Range("S71").Value = 2500
$ws.Range("T71").Value = 3

# New row 72: Segunda, week of 2023-08-03 (serial 45141)
$ws.Range("A72").Value = 1
$ws.Range("B72").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C72").Value = "Arica y Parinacota"
$ws.Range("D72").Value = 45141
$ws.Range("E72").Value = 15
$ws.Range("F72").Value = "Fruta"
$ws.Range("G72").Value = 100101
$ws.Range("H72").Value = "Berries"
$ws.Range("I72").Value = 100112025
$ws.Range("J72").Value = "Frutilla"
$ws.Range("K72").Value = "Sin especificar"
$ws.Range("L72").Value = "Segunda"
$ws.Range("M72").Value = 80
$ws.Range("N72").Value = 5000
$ws.Range("O72").Value = 6000
$ws.Range("P72").Value = 5500
$ws.Range("Q72").Value = "$/bandeja 3 kilos"
$ws.Range("R72").Value = "Región de Arica y Parinacota"
$ws.Range("S72").Value = 1833
$ws.Range("T72").Value = 3

# New row 73: Tercera, week of 2023-08-03 (serial 45141)
$ws.Range("A73").Value = 1
$ws.Range("B73").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C73").Value = "Arica y Parinacota"
$ws.Range("D73").Value = 45141
$ws.Range("E73").Value = 15
$ws.Range("F73").Value = "Fruta"
$ws.Range("G73").Value = 100101
$ws.Range("H73").Value = "Berries"
$ws.Range("I73").Value = 100112025
$ws.Range("J73").Value = "Frutilla"
$ws.Range("K73").Value = "Sin especificar"
$ws.Range("L73").Value = "Tercera"
$ws.Range("M73").Value = 100
$ws.Range("N73").Value = 3000
$ws.Range("O73").Value = 4000
$ws.Range("P73").Value = 3500
$ws.Range("Q73").Value = "$/bandeja 3 kilos"
$ws.Range("R73").Value = "Región de Arica y Parinacota"
$ws.Range("S73").Value = 1167
$ws.Range("T73").Value = 3
